$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# Replace the run text (also merges the two runs / drops the trailing
# space-only run, keeping the first run's character formatting).
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "**ID__AFFARS_AF_PGI_5345_103_70_90__ID**"

# Update the paragraph's left indent (225 twips == 11.25 points; Word's
# object model reports/accepts LeftIndent in points).
$p.Format.LeftIndent = 11.25

# Add the paragraph border (space-only border on all four sides, no
# visible line) matching the style already used elsewhere in the doc.
$borders = $p.Borders
$top = $borders.Item(-1)
$left = $borders.Item(-2)
$bottom = $borders.Item(-3)
$right = $borders.Item(-4)
$top.LineStyle = 0
$left.LineStyle = 0
$bottom.LineStyle = 0
$right.LineStyle = 0
$borders.DistanceFromTop = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromRight = 5
